$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for Wins / Losses / Ties (columns AD, AE, AF)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the bold/centered/bordered header style used by the rest of row 1
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Season record values for every data row (2-48)
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 83
    $ws.Cells.Item($r, 31).Value = 79
    $ws.Cells.Item($r, 32).Value = 0
}
